# Update the cached "date last saved" field text on the slide master and
# every slide layout (PowerPoint caches the rendered text of the
# datetimeFigureOut field in <a:t>, it does not recompute it here).
$p = $ppt.ActivePresentation

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = "06/08/2024"
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "06/08/2024"
        }
    }
}

# Slide 1: reposition the background picture and the "Python  LIB" caption,
# and update the caption text.
$s = $p.Slides.Item(1)

$pic = $s.Shapes.Item(1)
$pic.Top = 78.363622

$tb = $s.Shapes.Item(2)
$tb.Left = 298.4545
$tb.Top = 320.5454
$tb.Width = 545.909213
$tb.Height = 208.4156
$tb.TextFrame.TextRange.Text = "Python  LIB"
